$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.435.97'
$ws.Range('E2').Value = '  +1.90%  '
$ws.Range('D3').Value = '1.860.60'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.31'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4772'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3805'
$ws.Range('E8').Value = '  +3.81%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07314'
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9313'
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.73'
$ws.Range('E11').Value = '  +5.01%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07791'
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').Value = '1.873.10'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.445'
$ws.Range('E14').Value = '  +1.98%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.547'
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '90.13'
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008815'
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D20').Value = '27.489.28'
$ws.Range('E20').Value = '  +1.71%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.63'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.095'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.946'
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '154.96'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '18.45'
$ws.Range('E26').Value = '  +1.63%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.007'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '115.45'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.946'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.08895'
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.325'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.197'
$ws.Range('E32').Value = '  +2.02%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7533'
$ws.Range('E33').Value = '  +1.91%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.591'
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.696'
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.124'
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02040'
$ws.Range('E37').Value = '  +4.16%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5573'
$ws.Range('E38').Value = '  +6.55%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05278'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.990'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.036'
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.591'
$ws.Range('E42').Value = '  +4.05%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1520'
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4866'
$ws.Range('E44').Value = '  +2.94%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.62'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.662'
$ws.Range('E47').Value = '  +3.75%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '103.03'
$ws.Range('E48').Value = '  +1.41%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '67.44'
$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06089'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9136'
$ws.Range('E51').Value = '  +2.95%  '
